$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, using the same bold/centered/bordered style as the other headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for columns I (I0) and J (IF)
$values = @(
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(5, 6),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(1, 3),
    @(7, 7),
    @(6, 6),
    @(4, 5),
    @(5, 6),
    @(1, 2),
    @(3, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
